$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = 13.411
$ws.Range("E10").Value = 12.671
$ws.Range("E12").Value = 12.966
$ws.Range("E18").Value = 12.966
